$d = $word.ActiveDocument

# 1. Replace the literal merge field placeholder text *Profession* with guillemet-style «Profession»
$d.Content.Find.Execute("*Profession*", $true, $false, $false, $false, $false, $true, 1, $false, "«Profession»", 2)

# 2. Replace the plain hyphen in "Test hyphen - test." with an en dash
$d.Content.Find.Execute("Test hyphen - test.", $true, $false, $false, $false, $false, $true, 1, $false, "Test hyphen – test.", 2)

# 3. Replace straight double quotes with smart (curly) double quotes
$d.Content.Find.Execute('Test "Smart double quotes".', $true, $false, $false, $false, $false, $true, 1, $false, "Test “Smart double quotes”.", 2)

# 4. Replace straight single quotes with smart (curly) single quotes
$d.Content.Find.Execute("Test 'smart single quotes'.", $true, $false, $false, $false, $false, $true, 1, $false, "Test ‘smart single quotes’.", 2)
